# Populate the "castor" (beaver) row (row 3) of the metadata table with
# the new "Suivi du castor d'Europe" monitoring-network description, mirroring
# the level of detail already present for the other species rows (e.g. row 6,
# "becasse", and row 18, "onde").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 'Suivi du castor d''Europe'
$ws.Range("C3").Value = 'data-raw/logo_castor.jpg'
$ws.Range("D3").Value = 'Le Castor d''Europe est un mammifère semi-aquatique, et l''un des plus grand rongeur de la planète.
Cette espèce protégée est une "ingénieure des écosystèmes" au travers des différents aménagements qu''elle réalise dans les cours d''eau.'
$ws.Range("E3").Value = 'Les objectifs du réseau ont été fixés par le ministère en charge de l’écologie :
- assurer le suivi de la colonisation du castor sur le réseau hydrographique français ;
- accompagnement sur la question des dommages;
- vigilance sur l''arrivée du castor canadien'
$ws.Range("F3").Value = 'Raportage Directive européenne Habitat Faune Flore
Régulation des dispositifs de piègeage près des cours d''eau'
$ws.Range("G3").Value = 91
$ws.Range("I3").Value = 'Prospection de linéaires de cours d''eau à la recherche d''indices de présence'
$ws.Range("J3").Value = '11,12,1,2,3'
$ws.Range("L3").Value = 'Animation nationale:
Paul Hurel
Suivi scientifique:
Yoann Bressan
Animation régionale:
Cédric Mondy'
$ws.Range("M3").Value = 'ENS
Syndicats de rivière'
$ws.Range("X3").Value = 'ROE'
$ws.Range("Z3").Value = 'https://carmen.carmencarto.fr/38/Castor.map'
$ws.Range("AA3").Value = 'texte:Le réseau castor;lien:https://professionnels.ofb.fr/fr/reseau-castor'
$ws.Range("AB3").Value = 'texte:Fiche espèce;lien:https://professionnels.ofb.fr/fr/doc-fiches-especes/castor-deurope-castor-fiber'

# Row 3 now holds as much wrapped text as row 6 ("becasse"); match that row's
# height (288pt) rather than leaving the sheet's auto-computed height.
$ws.Rows.Item(3).RowHeight = 288

# Column E ("objectif") needs to be wide enough to show the new beaver-network
# objectives text; widen it to match the author's column width (~23.4 chars).
$ws.Columns.Item(5).ColumnWidth = 22.67

# Scroll the view down to the newly-filled-in row and reselect the objectif
# cell for that row, matching the author's final cursor position.
$ws.Range("E7").Select()

